# edit.ps1 - applies the authored changes described by the diff:
#  * Slide 10 (TextBox "TextBox 6"): wording tweak about identity/new user.
#  * Slide 11 (TextBox "TextBox 5"): "Update database" -> "update database".
#  * Slide 13 (TextBox "TextBox 3"): inserts "ContactOperations class under "
#    before "ContactManager.Authorization" and grows the box (autofit).
#  * Slide 14 (TextBox "TextBox 7"): "Running the App" -> longer caption and
#    the box is re-centered/re-sized (autofit, wrap="none").
#  * Slide 15 (Content placeholder): "Use cases" -> "Scenarios".

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# Slide 10: "...IdentityBuilder class. This is because when an identity is
# created it may belong to one or more roles." ->
# "...IdentityBuilder class. This is because when a new user is created, it
# may belong to one or more roles."
# ---------------------------------------------------------------------------
$s10 = $p.Slides.Item(10)
$shp10 = $s10.Shapes.Item(3)
$tr10 = $shp10.TextFrame.TextRange
$run10 = $tr10.Find(" class. This is because when an identity is created it may belong to one or more roles.")
$run10.Text = " class. This is because when a new user is created, it may belong to one or more roles."

# ---------------------------------------------------------------------------
# Slide 11: "Update database" -> "update database"
# ---------------------------------------------------------------------------
$s11 = $p.Slides.Item(11)
$shp11 = $s11.Shapes.Item(3)
$tr11 = $shp11.TextFrame.TextRange
$run11 = $tr11.Find("Update database")
$run11.Text = "update database"

# ---------------------------------------------------------------------------
# Slide 13: insert "ContactOperations class under " right before
# "ContactManager.Authorization", and grow the (autofit) textbox to match
# the extra wrapped line of text.
# ---------------------------------------------------------------------------
$s13 = $p.Slides.Item(13)
$shp13 = $s13.Shapes.Item(1)
$tr13 = $shp13.TextFrame.TextRange
# Anchor on the *entire* original run (not just a trailing substring) so
# the insertion point sits at the existing run boundary and the untouched
# run text isn't itself re-split.
$anchor13 = $tr13.Find(". The set of allowable operations in the Contact Manager App are configured in ")
$insertStart13 = $anchor13.Start + $anchor13.Text.Length
$null = $anchor13.InsertAfter("ContactOperations class under ")
# Split the freshly-inserted text into its own two runs (matching how
# PowerPoint breaks runs as formatting/spellcheck state changes while
# typing) by touching a formatting property on each sub-range with its
# own (unchanged) value.
$newRun13a = $tr13.Characters($insertStart13, 17)
$newRun13a.Font.Size = $newRun13a.Font.Size
$newRun13b = $tr13.Characters($insertStart13 + 17, 13)
$newRun13b.Font.Size = $newRun13b.Font.Size
$shp13.Height = 312.62346456692916

# ---------------------------------------------------------------------------
# Slide 14: "Running the App" -> "Running the App with server listening on
# port  53020"; the caption textbox uses wrap="none" + autofit so PowerPoint
# recenters/regrows it around its original center point once the text is
# longer.
# ---------------------------------------------------------------------------
$s14 = $p.Slides.Item(14)
$shp14 = $s14.Shapes.Item(5)
$shp14.TextFrame.TextRange.Text = "Running the App with server listening on port  53020"
$shp14.Left = 542.7543307086614
$shp14.Top = 499.20065
$shp14.Width = 411.5804724409449
$shp14.Height = 29.081259842519685

# ---------------------------------------------------------------------------
# Slide 15: "Use cases" -> "Scenarios"
# ---------------------------------------------------------------------------
$s15 = $p.Slides.Item(15)
$shp15 = $s15.Shapes.Item(1)
$tr15 = $shp15.TextFrame.TextRange
$run15 = $tr15.Find("Use cases")
$run15.Text = "Scenarios"
